$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 198 (row 197 "Segunda" stays put; everything
# from old row 198 onward shifts down by one, including the dangling last
# row which becomes the new row 276).
$ws.Rows.Item(198).Insert()

# Populate the freshly inserted row 198 with the new record.
$ws.Cells.Item(198, 1).Value = 4
$ws.Cells.Item(198, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(198, 3).Value = "Los Lagos"
$ws.Cells.Item(198, 4).Value = 44876
$ws.Cells.Item(198, 5).Value = 10
$ws.Cells.Item(198, 6).Value = 100112024
$ws.Cells.Item(198, 7).Value = "Choclo"
$ws.Cells.Item(198, 8).Value = "Dulce o Americano"
$ws.Cells.Item(198, 9).Value = "Primera"
$ws.Cells.Item(198, 10).Value = 180
$ws.Cells.Item(198, 11).Value = 30000
$ws.Cells.Item(198, 12).Value = 30000
$ws.Cells.Item(198, 13).Value = 30000
$ws.Cells.Item(198, 14).Value = "$/malla 70 unidades"
$ws.Cells.Item(198, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(198, 16).Value = 429
$ws.Cells.Item(198, 17).Value = 70
$ws.Cells.Item(198, 18).Value = "Hortaliza"
